$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "1.002", "48.20") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.271.70'
$ws.Range("E2").Value = '  +15.02%  '
$ws.Range("D3").Value = '1.679.76'
$ws.Range("E3").Value = '  +9.35%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -1.01%  '
$ws.Range("D5").Value = '307.72'
$ws.Range("E5").Value = '  +9.14%  '
$ws.Range("D6").Value = '0.9976'
$ws.Range("E6").Value = '  +3.40%  '
$ws.Range("D7").Value = '0.3728'
$ws.Range("E7").Value = '  +2.86%  '
$ws.Range("D8").Value = '0.3437'
$ws.Range("E8").Value = '  +8.07%  '
$ws.Range("D9").Value = '48.20'
$ws.Range("E9").Value = '  +18.66%  '
$ws.Range("D10").Value = '1.184'
$ws.Range("E10").Value = '  +8.35%  '
$ws.Range("D11").Value = '0.07299'
$ws.Range("E11").Value = '  +7.34%  '
$ws.Range("D12").Value = '0.9984'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").Value = '20.63'
$ws.Range("E13").Value = '  +10.41%  '
$ws.Range("D14").Value = '6.090'
$ws.Range("E14").Value = '  +7.54%  '
$ws.Range("D15").Value = '6.742'
$ws.Range("E15").Value = '  +6.22%  '
$ws.Range("D16").Value = '1.674.85'
$ws.Range("E16").Value = '  +9.88%  '
$ws.Range("D17").Value = '0.00001111'
$ws.Range("E17").Value = '  +6.32%  '
$ws.Range("D18").Value = '0.9976'
$ws.Range("E18").Value = '  +3.44%  '
$ws.Range("D19").Value = '0.06722'
$ws.Range("E19").Value = '  +10.29%  '
$ws.Range("D20").Value = '81.69'
$ws.Range("E20").Value = '  +13.24%  '
$ws.Range("D21").Value = '16.48'
$ws.Range("E21").Value = '  +10.11%  '
$ws.Range("D22").Value = '6.125'
$ws.Range("E22").Value = '  +7.50%  '
$ws.Range("D23").Value = '12.03'
$ws.Range("E23").Value = '  +5.98%  '
$ws.Range("D24").Value = '24.205.80'
$ws.Range("E24").Value = '  +14.52%  '
$ws.Range("D25").Value = '2.402'
$ws.Range("E25").Value = '  +3.15%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.675'
$ws.Range("E26").Value = '  +20.89%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '3.359'
$ws.Range("E27").Value = '  -8.99%  '
$ws.Range("D28").Value = '152.18'
$ws.Range("E28").Value = '  +2.66%  '
$ws.Range("E29").Value = '  +10.52%  '
$ws.Range("D30").Value = '1.858.55'
$ws.Range("E30").Value = '  +9.75%  '
$ws.Range("D31").Value = '126.95'
$ws.Range("E31").Value = '  +7.22%  '
$ws.Range("D32").Value = '6.420'
$ws.Range("E32").Value = '  +23.94%  '
$ws.Range("D33").Value = '4.043'
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").Value = '0.9909'
$ws.Range("E34").Value = '  +16.50%  '
$ws.Range("D35").Value = '1.736'
$ws.Range("E35").Value = '  +15.20%  '
$ws.Range("D36").Value = '0.08445'
$ws.Range("E36").Value = '  +5.62%  '
$ws.Range("D37").Value = '12.45'
$ws.Range("E37").Value = '  +17.01%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '5.366'
$ws.Range("E38").Value = '  +8.16%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.06414'
$ws.Range("E39").Value = '  +9.44%  '
$ws.Range("D40").Value = '8.885'
$ws.Range("E40").Value = '  +15.10%  '
$ws.Range("D41").Value = '1.290'
$ws.Range("E41").Value = '  +7.45%  '
$ws.Range("D42").Value = '0.02342'
$ws.Range("E42").Value = '  +11.40%  '
$ws.Range("D43").Value = '0.2116'
$ws.Range("E43").Value = '  +10.69%  '
$ws.Range("D44").Value = '0.6158'
$ws.Range("E44").Value = '  +13.37%  '
$ws.Range("D45").Value = '0.9969'
$ws.Range("E45").Value = '  +3.37%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '3.804'
$ws.Range("E46").Value = '  +6.19%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '13.20'
$ws.Range("E47").Value = '  +5.07%  '
$ws.Range("D48").Value = '0.5964'
$ws.Range("E48").Value = '  +9.61%  '
$ws.Range("D49").Value = '127.88'
$ws.Range("E49").Value = '  +5.27%  '
$ws.Range("D50").Value = '2.022'
$ws.Range("E50").Value = '  +8.22%  '
$ws.Range("D51").Value = '0.07153'
$ws.Range("E51").Value = '  +8.89%  '
